$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column B (old B:E shift right to E:H)
$ws.Range("B1:D1").EntireColumn.Insert()

# Fill the newly inserted (blank) data columns with the "UN" placeholder used
# throughout the watchlist for not-yet-updated cells
$ws.Range("B2:D27").Value = "UN"

# New header values for the freshly inserted columns (most-recent-first date columns)
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# New analyst group appended at the bottom of the watchlist
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"

# New latest-date header column set last
$ws.Range("B1").Value = "Jun_27"
